# preparation publication 0.2.0
# - bump Version 0.1.1 -> 0.2.0
# - bump Date to the new publication timestamp
# - add a new "Jurisdiction" / "iso:code:3166:FR" metadata row
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Version value (row 3, column B)
$ws.Cells.Item(3, 2).Value = "0.2.0"

# Update the Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2023-10-20T08:59:58+00:00"

# Insert a new row after "Contact" (row 10) for the Jurisdiction property,
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting of the other data rows (thin border, top-aligned, wrapped text)
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Set the actual values for the new Jurisdiction row
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = "iso:code:3166:FR"
